# Fruta / hortaliza, semanal
#
# A new week's worth of price observations (2 rows: Primera + Segunda
# quality, date 2022-09-05 = serial 44809) is prepended to the
# date-ordered block of "Kiwi" rows that starts at row 266. Every
# existing data row from 266 downward shifts down by two rows
# (266->268 ... 358->360), and the sheet's used range grows from
# A1:T358 to A1:T360.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 266/267, pushing the existing rows (and
# everything below them) down by two.
$ws.Range("A266:A267").EntireRow.Insert()

# --- New row 266 -----------------------------------------------------
$ws.Range("A266").Value = 4
$ws.Range("B266").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C266").Value = "Los Lagos"
$ws.Range("D266").Value = 44809
$ws.Range("E266").Value = 10
$ws.Range("F266").Value = "Fruta"
$ws.Range("G266").Value = 100101
$ws.Range("H266").Value = "Berries"
$ws.Range("I266").Value = 100101007
$ws.Range("J266").Value = "Kiwi"
$ws.Range("K266").Value = "Hayward"
$ws.Range("L266").Value = "Primera"
$ws.Range("M266").Value = 100
$ws.Range("N266").Value = 14000
$ws.Range("O266").Value = 14000
$ws.Range("P266").Value = 14000
$ws.Range("Q266").Value = "`$/caja 15 kilos"
$ws.Range("R266").Value = "Región de O'Higgins"
$ws.Range("S266").Value = 933
$ws.Range("T266").Value = 15

# --- New row 267 -----------------------------------------------------
$ws.Range("A267").Value = 4
$ws.Range("B267").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C267").Value = "Los Lagos"
$ws.Range("D267").Value = 44809
$ws.Range("E267").Value = 10
$ws.Range("F267").Value = "Fruta"
$ws.Range("G267").Value = 100101
$ws.Range("H267").Value = "Berries"
$ws.Range("I267").Value = 100101007
$ws.Range("J267").Value = "Kiwi"
$ws.Range("K267").Value = "Hayward"
$ws.Range("L267").Value = "Segunda"
$ws.Range("M267").Value = 100
$ws.Range("N267").Value = 12500
$ws.Range("O267").Value = 12500
$ws.Range("P267").Value = 12500
$ws.Range("Q267").Value = "`$/caja 15 kilos"
$ws.Range("R267").Value = "Región de O'Higgins"
$ws.Range("S267").Value = 833
$ws.Range("T267").Value = 15
